# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / handoff / handback timestamps
# that get refreshed each time the handback report is (re)generated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G4 and de-de!H4 both show the "Latest HO Xliff Generate Date" /
# "Correspond Handoff Datetime" for d8e05e25-...md (they share the same text).
$wsOverview.Range("G4").Value = "2016-08-29 18:49:41"
$wsDeDe.Range("H4").Value = "2016-08-29 18:49:41"

# zh-cn!H4 - Correspond Handoff Datetime for d8e05e25-...md
$wsZhCn.Range("H4").Value = "2016-08-29 18:49:36"

# zh-cn!K4 - Correspond Handback DateTime for d8e05e25-...md
$wsZhCn.Range("K4").Value = "2016-08-29 18:50:13"

# de-de!K4 - Correspond Handback DateTime for d8e05e25-...md
$wsDeDe.Range("K4").Value = "2016-08-29 18:50:25"
